# Add a new 4th slide ("Github") using the "Section Header" layout
# (same layout already used as CustomLayout 3 / slideLayout3.xml, type="secHead"),
# with a title "Github" and a body placeholder holding hyperlinked text
# pointing at the team's GitHub repo.

$p = $ppt.ActivePresentation

# "Section Header" is CustomLayout index 3 on the (single) slide master -
# it provides exactly the two placeholders we need: title + body(idx=1).
$master = $p.SlideMaster
$sectionHeaderLayout = $master.CustomLayouts.Item(3)

# Append the new slide at the end of the deck (becomes slide 4).
$s = $p.Slides.AddSlide($p.Slides.Count + 1, $sectionHeaderLayout)

# Title placeholder -> "Github"
$titleShape = $s.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "Github"

# Body placeholder -> hyperlinked repo text
$bodyShape = $s.Shapes.Item(2)
$bodyRange = $bodyShape.TextFrame.TextRange
$bodyRange.Text = "GitHub - KevTuco/STDISCM-P2"
$bodyRange.ActionSettings.Item(1).Hyperlink.Address = "https://github.com/KevTuco/STDISCM-P2"
